# Applies the commit's edits to the "Hoja1" worksheet:
#  - Row 17 (document/name pair) is changed to 1128052019 / JENNIFER TATIANA HENAO QUICENO
#  - Row 18 (document/name pair) is changed to 1002076986 / ROBERT CARABALLO PADILLA
#  - G18 (Salario Basico for row 18) changes from 1250000 to 908526

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the worker identity on rows 17 and 18
$ws.Range("C17").Value = "1128052019"
$ws.Range("D17").Value = "JENNIFER TATIANA HENAO QUICENO"

$ws.Range("C18").Value = "1002076986"
$ws.Range("D18").Value = "ROBERT CARABALLO PADILLA"

# Update the "Salario Basico" value for row 18
$ws.Range("G18").Value = 908526
